$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new note text to C11
$ws.Range("C11").Value = "automate the fill in function, like def_(): to iterate through all features, and all dependent features"

# Update the active cell selection to A11
$ws.Range("A11").Select()

$wb.Save()
